# Fruta / hortaliza, semanal
# Insert two new weekly data rows (203:204) into the Cilantro log, pushing
# the existing rows 203:220 down to 205:222, then populate the two new
# rows with the latest observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 203 (shifts old 203-220 -> 205-222)
$ws.Rows("203:204").Insert()

# New row 203
$ws.Range("A203").Value2 = 10
$ws.Range("B203").Value2 = 'Vega Modelo de Temuco'
$ws.Range("C203").Value2 = 'La Araucanía'
$ws.Range("D203").Value2 = 44461
$ws.Range("E203").Value2 = 9
$ws.Range("F203").Value2 = 100112040
$ws.Range("G203").Value2 = 'Cilantro'
$ws.Range("H203").Value2 = 'Sin especificar'
$ws.Range("I203").Value2 = 'Primera'
$ws.Range("J203").Value2 = 25
$ws.Range("K203").Value2 = 6000
$ws.Range("L203").Value2 = 6000
$ws.Range("M203").Value2 = 6000
$ws.Range("N203").Value2 = '$/docena de atados (2 kilos)'
$ws.Range("O203").Value2 = 'Provincia de Cautín'
$ws.Range("P203").Value2 = 3000
$ws.Range("Q203").Value2 = 2
$ws.Range("R203").Value2 = 'Hortaliza'

# New row 204
$ws.Range("A204").Value2 = 10
$ws.Range("B204").Value2 = 'Vega Modelo de Temuco'
$ws.Range("C204").Value2 = 'La Araucanía'
$ws.Range("D204").Value2 = 44461
$ws.Range("E204").Value2 = 9
$ws.Range("F204").Value2 = 100112040
$ws.Range("G204").Value2 = 'Cilantro'
$ws.Range("H204").Value2 = 'Sin especificar'
$ws.Range("I204").Value2 = 'Primera'
$ws.Range("J204").Value2 = 40
$ws.Range("K204").Value2 = 4000
$ws.Range("L204").Value2 = 4000
$ws.Range("M204").Value2 = 4000
$ws.Range("N204").Value2 = '$/docena de atados (2 kilos)'
$ws.Range("O204").Value2 = 'Región Metropolitana'
$ws.Range("P204").Value2 = 2000
$ws.Range("Q204").Value2 = 2
$ws.Range("R204").Value2 = 'Hortaliza'
